$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 198, shifting existing rows 198-236 down to 199-237.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly record.
$ws.Range("A198").Value = 5
$ws.Range("B198").Value = "Macroferia Regional de Talca"
$ws.Range("C198").Value = "Maule"
$ws.Range("D198").Value = 44694
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 100112008
$ws.Range("G198").Value = "Coliflor"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 3000
$ws.Range("K198").Value = 1000
$ws.Range("L198").Value = 1000
$ws.Range("M198").Value = 1000
$ws.Range("N198").Value = "$/unidad"
$ws.Range("O198").Value = "Región del Maule"
$ws.Range("P198").Value = 1000
$ws.Range("Q198").Value = 1
$ws.Range("R198").Value = "Hortaliza"
